$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Cells.Item(2, 4).Value = 212
$ws.Cells.Item(2, 5).Value = 247
$ws.Cells.Item(2, 6).Value = 278
$ws.Cells.Item(2, 7).Value = 326

$ws.Cells.Item(3, 4).Value = 189
$ws.Cells.Item(3, 5).Value = 224
$ws.Cells.Item(3, 6).Value = 256

$ws.Cells.Item(4, 4).Value = 194
$ws.Cells.Item(4, 5).Value = 229
$ws.Cells.Item(4, 6).Value = 263
$ws.Cells.Item(4, 7).Value = 315

$ws.Cells.Item(5, 4).Value = 200
$ws.Cells.Item(5, 5).Value = 237
$ws.Cells.Item(5, 6).Value = 272
$ws.Cells.Item(5, 7).Value = 326

$ws.Cells.Item(6, 4).Value = 208
$ws.Cells.Item(6, 5).Value = 248
$ws.Cells.Item(6, 6).Value = 288
$ws.Cells.Item(6, 7).Value = 349

$ws.Cells.Item(7, 4).Value = 209
$ws.Cells.Item(7, 5).Value = 247
$ws.Cells.Item(7, 6).Value = 284
$ws.Cells.Item(7, 7).Value = 342

$ws.Cells.Item(8, 4).Value = 215
$ws.Cells.Item(8, 5).Value = 258
$ws.Cells.Item(8, 6).Value = 301
$ws.Cells.Item(8, 7).Value = 368

$ws.Cells.Item(9, 4).Value = 217
$ws.Cells.Item(9, 5).Value = 259
$ws.Cells.Item(9, 6).Value = 302
$ws.Cells.Item(9, 7).Value = 368

$ws.Cells.Item(10, 4).Value = 207
$ws.Cells.Item(10, 5).Value = 246
$ws.Cells.Item(10, 6).Value = 286
$ws.Cells.Item(10, 7).Value = 346

$ws.Cells.Item(11, 4).Value = 213
$ws.Cells.Item(11, 5).Value = 255
$ws.Cells.Item(11, 6).Value = 298
$ws.Cells.Item(11, 7).Value = 365

$ws.Cells.Item(12, 4).Value = 216
$ws.Cells.Item(12, 5).Value = 259
$ws.Cells.Item(12, 6).Value = 304
$ws.Cells.Item(12, 7).Value = 372

$ws.Cells.Item(13, 4).Value = 216
$ws.Cells.Item(13, 5).Value = 262
$ws.Cells.Item(13, 6).Value = 314
$ws.Cells.Item(13, 7).Value = 396

$ws.Cells.Item(14, 4).Value = 214
$ws.Cells.Item(14, 5).Value = 259
$ws.Cells.Item(14, 6).Value = 310
$ws.Cells.Item(14, 7).Value = 390

$ws.Cells.Item(15, 4).Value = 210
$ws.Cells.Item(15, 5).Value = 255
$ws.Cells.Item(15, 6).Value = 306
$ws.Cells.Item(15, 7).Value = 388

$ws.Cells.Item(16, 4).Value = 205
$ws.Cells.Item(16, 5).Value = 249
$ws.Cells.Item(16, 6).Value = 299
$ws.Cells.Item(16, 7).Value = 379

$ws.Cells.Item(17, 4).Value = 199
$ws.Cells.Item(17, 5).Value = 242
$ws.Cells.Item(17, 6).Value = 294
$ws.Cells.Item(17, 7).Value = 376

